$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.03689613626135
$ws.Range("C2").Value = 11.52911636749117
$ws.Range("D2").Value = 5.972085103660648
$ws.Range("E2").Value = 16.50818836825585
$ws.Range("G2").Value = 39.22750668874065
$ws.Range("H2").Value = 15.56785156535916
$ws.Range("N2").Value = 16.70484267705587
$ws.Range("B3").Value = 16.2738382835976
$ws.Range("C3").Value = 10.76091934065724
$ws.Range("D3").Value = 5.852429783991138
$ws.Range("E3").Value = 15.56474104257025
$ws.Range("G3").Value = 38.21224104383481
$ws.Range("H3").Value = 15.50067658952727
$ws.Range("N3").Value = 16.78284731500795
$ws.Range("B4").Value = 15.79170575756254
$ws.Range("C4").Value = 10.26013374642917
$ws.Range("D4").Value = 5.779761241885162
$ws.Range("E4").Value = 14.96142054577128
$ws.Range("G4").Value = 37.59329657275025
$ws.Range("H4").Value = 15.46437675690437
$ws.Range("N4").Value = 16.8327148442062
$ws.Range("B5").Value = 15.59212323495235
$ws.Range("C5").Value = 10.04871409016531
$ws.Range("D5").Value = 5.750393338361255
$ws.Range("E5").Value = 14.70978760093013
$ws.Range("G5").Value = 37.3426154512547
$ws.Range("H5").Value = 15.45083247941425
$ws.Range("N5").Value = 16.85353389926881
$ws.Range("B6").Value = 15.55880485317855
$ws.Range("C6").Value = 10.01316361757609
$ws.Range("D6").Value = 5.745532908848067
$ws.Range("E6").Value = 14.66766429965372
$ws.Range("G6").Value = 37.30109598732123
$ws.Range("H6").Value = 15.44865894427571
$ws.Range("N6").Value = 16.85702099211328
$ws.Range("B7").Value = 15.78902626462461
$ws.Range("C7").Value = 10.25731224275814
$ws.Range("D7").Value = 5.779364125534194
$ws.Range("E7").Value = 14.95804992236081
$ws.Range("G7").Value = 37.58990897530204
$ws.Range("H7").Value = 15.46418903556892
$ws.Range("N7").Value = 16.83299360013762
$ws.Range("B8").Value = 16.77679444426207
$ws.Range("C8").Value = 11.27026826622946
$ws.Range("D8").Value = 5.930688098030909
$ws.Range("E8").Value = 16.18802346645868
$ws.Range("G8").Value = 38.87677746354422
$ws.Range("H8").Value = 15.54366640994965
$ws.Range("N8").Value = 16.73133031967497
$ws.Range("B9").Value = 18.59402542031144
$ws.Range("C9").Value = 13.02685772686901
$ws.Range("D9").Value = 6.231876580794786
$ws.Range("E9").Value = 18.49463885423302
$ws.Range("G9").Value = 41.41724141192646
$ws.Range("H9").Value = 15.73847778379691
$ws.Range("N9").Value = 16.54754084595467
$ws.Range("B10").Value = 19.84287966778541
$ws.Range("C10").Value = 14.17892224268069
$ws.Range("D10").Value = 6.45337621046393
$ws.Range("E10").Value = 20.15921214456702
$ws.Range("G10").Value = 43.27078541270711
$ws.Range("H10").Value = 15.90484708400396
$ws.Range("N10").Value = 16.42189194404978
$ws.Range("B11").Value = 20.39012444712388
$ws.Range("C11").Value = 14.67329237863586
$ws.Range("D11").Value = 6.553688660743695
$ws.Range("E11").Value = 20.87494270755589
$ws.Range("G11").Value = 44.10670569162296
$ws.Range("H11").Value = 15.98543370933745
$ws.Range("N11").Value = 16.36674415396013
$ws.Range("B12").Value = 20.59420752125033
$ws.Range("C12").Value = 14.85625184536864
$ws.Range("D12").Value = 6.591570580804554
$ws.Range("E12").Value = 21.14004664248803
$ws.Range("G12").Value = 44.42186088355119
$ws.Range("H12").Value = 16.01664006878525
$ws.Range("N12").Value = 16.34614848842054
$ws.Range("B13").Value = 20.5503966092342
$ws.Range("C13").Value = 14.81703678321944
$ws.Range("D13").Value = 6.583417276250433
$ws.Range("E13").Value = 21.08321478565921
$ws.Range("G13").Value = 44.35405359027966
$ws.Range("H13").Value = 16.00988879476456
$ws.Range("N13").Value = 16.35057136589021
$ws.Range("B14").Value = 20.40697830267578
$ws.Range("C14").Value = 14.68842955606942
$ws.Range("D14").Value = 6.556807537710359
$ws.Range("E14").Value = 20.89687154559293
$ws.Range("G14").Value = 44.13266343152023
$ws.Range("H14").Value = 15.98798734067935
$ws.Range("N14").Value = 16.36504398068081
$ws.Range("B15").Value = 20.31871672239593
$ws.Range("C15").Value = 14.60910145449055
$ws.Range("D15").Value = 6.540493589483024
$ws.Range("E15").Value = 20.781960120273
$ws.Range("G15").Value = 43.99686480557786
$ws.Range("H15").Value = 15.9746614397717
$ws.Range("N15").Value = 16.37394628549293
$ws.Range("B16").Value = 19.80668302639473
$ws.Range("C16").Value = 14.14601824771769
$ws.Range("D16").Value = 6.446808083715887
$ws.Range("E16").Value = 20.11160586200481
$ws.Range("G16").Value = 43.21597976891211
$ws.Range("H16").Value = 15.89967797579741
$ws.Range("N16").Value = 16.42553630962437
$ws.Range("B17").Value = 19.48711112400809
$ws.Range("C17").Value = 13.85433974436203
$ws.Range("D17").Value = 6.389191204712596
$ws.Range("E17").Value = 19.6897657581402
$ws.Range("G17").Value = 42.73481742699281
$ws.Range("H17").Value = 15.85492380123183
$ws.Range("N17").Value = 16.45769898517346
$ws.Range("B18").Value = 19.30134499726179
$ws.Range("C18").Value = 13.68377610118912
$ws.Range("D18").Value = 6.356011499509481
$ws.Range("E18").Value = 19.44322673221755
$ws.Range("G18").Value = 42.45740063388718
$ws.Range("H18").Value = 15.82964502765886
$ws.Range("N18").Value = 16.47638743462628
$ws.Range("B19").Value = 19.23811654159995
$ws.Range("C19").Value = 13.62554489669379
$ws.Range("D19").Value = 6.344771794828846
$ws.Range("E19").Value = 19.35908050924478
$ws.Range("G19").Value = 42.36336880524333
$ws.Range("H19").Value = 15.82116596899526
$ws.Range("N19").Value = 16.48274758605832
$ws.Range("B20").Value = 19.52133375658785
$ws.Range("C20").Value = 13.88567881163303
$ws.Range("D20").Value = 6.39532905389707
$ws.Range("E20").Value = 19.73507551344093
$ws.Range("G20").Value = 42.78610937964697
$ws.Range("H20").Value = 15.85964018348843
$ws.Range("N20").Value = 16.45425562958226
$ws.Range("B21").Value = 20.44919018471665
$ws.Range("C21").Value = 14.72631970514639
$ws.Range("D21").Value = 6.564626597361384
$ws.Range("E21").Value = 20.95176563929716
$ws.Range("G21").Value = 44.1977314205115
$ws.Range("H21").Value = 15.99440173297304
$ws.Range("N21").Value = 16.36078522925632
$ws.Range("B22").Value = 21.03720497264508
$ws.Range("C22").Value = 15.25097912135801
$ws.Range("D22").Value = 6.674647938148601
$ws.Range("E22").Value = 21.71240977215843
$ws.Range("G22").Value = 45.11206667744559
$ws.Range("H22").Value = 16.08648962400661
$ws.Range("N22").Value = 16.30137262003969
$ws.Range("B23").Value = 20.7250957887573
$ws.Range("C23").Value = 14.97321490705268
$ws.Range("D23").Value = 6.615997162422136
$ws.Range("E23").Value = 21.30958667103472
$ws.Range("G23").Value = 44.62492991495954
$ws.Range("H23").Value = 16.03697879657852
$ws.Range("N23").Value = 16.33292941345728
$ws.Range("B24").Value = 19.50586804715138
$ws.Range("C24").Value = 13.87151937237503
$ws.Range("D24").Value = 6.392554300540505
$ws.Range("E24").Value = 19.71460348972188
$ws.Range("G24").Value = 42.76292272254998
$ws.Range("H24").Value = 15.85750650117246
$ws.Range("N24").Value = 16.45581175474329
$ws.Range("B25").Value = 18.11673610630323
$ws.Range("C25").Value = 12.57609131042948
$ws.Range("D25").Value = 6.150183572389892
$ws.Range("E25").Value = 17.84459267966255
$ws.Range("G25").Value = 40.73064510304653
$ws.Range("H25").Value = 15.68165067321743
$ws.Range("N25").Value = 16.59560530536929
